$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: re-labelled metadata "kind" strings for several columns -------
# C2: temporalidad column is now classified as a measure, not a dimension
$ws.Range("C2").Value = "iaest-measure:temporalidad"
# H2: sector-descripcion column is now classified as a measure, not a dimension
$ws.Range("H2").Value = "iaest-measure:sector-descripcion"
# M2: direccion-provincial-nombre column now gets its own curated measure label
$ws.Range("M2").Value = "iaest-measure:direccion-provincial-nombre"

# --- Row 3: "medida" (measure) replaces "dim" (dimension) for the columns
#     that were reclassified above (plus the columns that already shared the
#     same shared-string slot) ------------------------------------------------
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("H3").Value = "medida"
$ws.Range("L3").Value = "medida"
$ws.Range("M3").Value = "medida"

# --- Row 4: "xsd:int" replaces "skos:Concept"/"URI-Provincia" for the same
#     reclassified columns, since they are now plain integer measures -------
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("H4").Value = "xsd:int"
$ws.Range("L4").Value = "xsd:int"
$ws.Range("M4").Value = "xsd:int"

# --- Row 5: only the "ano" dimension keeps a mapping workbook now; the
#     temporalidad/sector-descripcion mapping files are no longer relevant ---
$ws.Range("C5").Clear()
$ws.Range("E5").Value = "mapping-ano.xlsx"
$ws.Range("H5").Clear()
